# Liga Finlandia 2025 - add the 2025-07-18 .. 2025-07-21 matchday results
# (rows 93-98) to Sheet1, matching the new data appended to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Each inner array lines up with $cols:
# Fecha, Local, Visita, Goles Local, Goles Visita, Fixture ID,
# Corners Local, Corners Visita, Amarillas Local, Amarillas Visita,
# Rojas Local, Rojas Visita, Goles 1T Local, Goles 1T Visita,
# Goles 2T Local, Goles 2T Visita, Posesion Local (%), Posesion Visita (%),
# Posesion Visita (%).1, Resultado
$data = @(
    @("2025-07-18", "Kooteepee",     "Inter Turku", 2, 2, 1342756, 4,  11, 3, 2, 0, 0, 0, 0, 2, 2, "42%", "58%", " ", "E"),
    @("2025-07-19", "VPS",           "KuPS",        1, 1, 1342757, 4,  1,  2, 2, 0, 1, 0, 0, 1, 1, "54%", "46%", " ", "E"),
    @("2025-07-20", "SJK",           "Gnistan",     3, 1, 1342759, 9,  1,  1, 1, 0, 0, 0, 0, 3, 1, "60%", "40%", " ", "L"),
    @("2025-07-20", "HJK helsinki",  "AC Oulu",     3, 1, 1342760, 4,  8,  1, 1, 0, 0, 0, 0, 3, 1, "50%", "50%", " ", "L"),
    @("2025-07-20", "FF Jaro",       "Ilves",       0, 1, 1342758, 1,  17, 3, 2, 1, 0, 0, 0, 0, 1, "36%", "64%", " ", "V"),
    @("2025-07-21", "Mariehamn",     "Haka",        1, 1, 1342761, 10, 5,  1, 6, 0, 0, 0, 0, 1, 1, "57%", "43%", " ", "E")
)

$startRow = 93
$endRow = $startRow + $data.Length - 1

# Columns A (Fecha), Q and R (Posesion %) hold values that would otherwise be
# auto-converted to dates / numbers by Excel, so force them to text first so
# the literal strings ("2025-07-18", "42%", ...) are preserved.
$ws.Range("A$startRow" + ":A$endRow").NumberFormat = "@"
$ws.Range("Q$startRow" + ":R$endRow").NumberFormat = "@"

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $startRow + $r
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ref = "$($cols[$c])$rowNum"
        $ws.Range($ref).Value = $rowValues[$c]
    }
}
